# Logged Week 16 and performed season sim from Week 17
# Update row 3 ("R") stats on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 432
$wsOff.Range("C3").Value = 294
$wsOff.Range("D3").Value = 113
$wsOff.Range("E3").Value = 52
$wsOff.Range("F3").Value = 7
$wsOff.Range("G3").Value = 4

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 407
$wsDef.Range("C3").Value = 274
$wsDef.Range("D3").Value = 78
$wsDef.Range("E3").Value = 32
